# Scheduled-runner update: refresh market-price-derived profit columns
# (H:currentAveragePrice, I/J:currentAveragePriceNQ/HQ, K/L:LevePriceNQ/HQ,
# M/N:LeveProfitNQ/HQ) on a handful of leve rows across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 2178.3704
$ws.Range("I98").Value = 1492.9231
$ws.Range("J98").Value = 20000
$ws.Range("K98").Value = 1492.9231
$ws.Range("L98").Value = 20000
$ws.Range("M98").Value = 5.076900000000023
$ws.Range("N98").Value = -22996

# Row 122
$ws.Range("H122").Value = 2178.3704
$ws.Range("I122").Value = 1492.9231
$ws.Range("J122").Value = 20000
$ws.Range("K122").Value = 4478.7693
$ws.Range("L122").Value = 60000
$ws.Range("M122").Value = -2028.7693
$ws.Range("N122").Value = -64900

$ws = $wb.Worksheets.Item("ARM")
# Row 97
$ws.Range("H97").Value = 861.9167
$ws.Range("I97").Value = 789.62067
$ws.Range("J97").Value = 1161.4286
$ws.Range("K97").Value = 789.62067
$ws.Range("L97").Value = 1161.4286
$ws.Range("M97").Value = -293.62067
$ws.Range("N97").Value = -2153.4286

# Row 132
$ws.Range("H132").Value = 2566977.8
$ws.Range("I132").Value = 1679.3077
$ws.Range("J132").Value = 7697575
$ws.Range("K132").Value = 5037.9231
$ws.Range("L132").Value = 23092725
$ws.Range("M132").Value = -2507.9231
$ws.Range("N132").Value = -23097785

$ws = $wb.Worksheets.Item("BSM")
# Row 10
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("N10").Value = 0

# Row 23
$ws.Range("H23").Value = 3574.3333
$ws.Range("I23").Value = 1500
$ws.Range("J23").Value = 4611.5
$ws.Range("K23").Value = 1500
$ws.Range("L23").Value = 4611.5
$ws.Range("M23").Value = -1217
$ws.Range("N23").Value = -5177.5

# Row 99
$ws.Range("H99").Value = 50001160
$ws.Range("I99").Value = 100000936
$ws.Range("J99").Value = 1388.7
$ws.Range("K99").Value = 100000936
$ws.Range("L99").Value = 1388.7
$ws.Range("M99").Value = -99999438
$ws.Range("N99").Value = -4384.7

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5040.717
$ws.Range("I31").Value = 1642.7222
$ws.Range("J31").Value = 12236.471
$ws.Range("K31").Value = 1642.7222
$ws.Range("L31").Value = 12236.471
$ws.Range("M31").Value = -1347.7222
$ws.Range("N31").Value = -12826.471

# Row 34
$ws.Range("H34").Value = 5040.717
$ws.Range("I34").Value = 1642.7222
$ws.Range("J34").Value = 12236.471
$ws.Range("K34").Value = 1642.7222
$ws.Range("L34").Value = 12236.471
$ws.Range("M34").Value = -1440.7222
$ws.Range("N34").Value = -12640.471

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 204899.67
$ws.Range("I5").Value = 408.14285
$ws.Range("J5").Value = 286696.28
$ws.Range("K5").Value = 1224.42855
$ws.Range("L5").Value = 860088.8400000001
$ws.Range("M5").Value = -1112.42855
$ws.Range("N5").Value = -860312.8400000001

# Row 23
$ws.Range("H23").Value = 5882470.5
$ws.Range("I23").Value = 16666750
$ws.Range("J23").Value = 136.45454
$ws.Range("K23").Value = 50000250
$ws.Range("L23").Value = 409.36362
$ws.Range("M23").Value = -50000015
$ws.Range("N23").Value = -879.3636200000001

# Row 81
$ws.Range("H81").Value = 2200
$ws.Range("I81").Value = 3800
$ws.Range("J81").Value = 1400
$ws.Range("K81").Value = 11400
$ws.Range("L81").Value = 4200
$ws.Range("M81").Value = -10277
$ws.Range("N81").Value = -6446

# Row 84
$ws.Range("H84").Value = 2200
$ws.Range("I84").Value = 3800
$ws.Range("J84").Value = 1400
$ws.Range("K84").Value = 34200
$ws.Range("L84").Value = 12600
$ws.Range("M84").Value = -28584
$ws.Range("N84").Value = -23832

# Row 97
$ws.Range("H97").Value = 4167021
$ws.Range("I97").Value = 4762252.5
$ws.Range("J97").Value = 400
$ws.Range("K97").Value = 14286757.5
$ws.Range("L97").Value = 1200
$ws.Range("M97").Value = -14286261.5
$ws.Range("N97").Value = -2192

# Row 131
$ws.Range("H131").Value = 2001031.5
$ws.Range("I131").Value = 6250683.5
$ws.Range("J131").Value = 1195.2941
$ws.Range("K131").Value = 18752050.5
$ws.Range("L131").Value = 3585.8823
$ws.Range("M131").Value = -18747010.5
$ws.Range("N131").Value = -13665.8823

# Row 132
$ws.Range("H132").Value = 1934.1666
$ws.Range("I132").Value = 2365.8572
$ws.Range("J132").Value = 1659.4546
$ws.Range("K132").Value = 21292.7148
$ws.Range("L132").Value = 14935.0914
$ws.Range("M132").Value = -18762.7148
$ws.Range("N132").Value = -19995.0914

# Row 135
$ws.Range("H135").Value = 204899.67
$ws.Range("I135").Value = 408.14285
$ws.Range("J135").Value = 286696.28
$ws.Range("K135").Value = 3673.28565
$ws.Range("L135").Value = 2580266.52
$ws.Range("M135").Value = -1138.28565
$ws.Range("N135").Value = -2585336.52

$ws = $wb.Worksheets.Item("GSM")
# Row 64
$ws.Range("H64").Value = 32521
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 32521
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 32521
$ws.Range("N64").Value = -33017

# Row 67
$ws.Range("H67").Value = 32521
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 32521
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 32521
$ws.Range("N67").Value = -34237

# Row 102
$ws.Range("H102").Value = 585888.5600000001
$ws.Range("I102").Value = 892913.9399999999
$ws.Range("J102").Value = 2540.3
$ws.Range("K102").Value = 892913.9399999999
$ws.Range("L102").Value = 2540.3
$ws.Range("M102").Value = -891291.9399999999
$ws.Range("N102").Value = -5784.3

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 4077578.2
$ws.Range("I122").Value = 4208327.5
$ws.Range("J122").Value = 3336666.8
$ws.Range("K122").Value = 12624982.5
$ws.Range("L122").Value = 10010000.4
$ws.Range("M122").Value = -12622532.5
$ws.Range("N122").Value = -10014900.4

# Row 134
$ws.Range("H134").Value = 41353.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 41353.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 41353.5
$ws.Range("N134").Value = -51493.5

# Row 141
$ws.Range("H141").Value = 57303.125
$ws.Range("I141").Value = 220000
$ws.Range("J141").Value = 34060.715
$ws.Range("K141").Value = 220000
$ws.Range("L141").Value = 34060.715
$ws.Range("M141").Value = -214820
$ws.Range("N141").Value = -44420.715

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1308.0834
$ws.Range("I122").Value = 1266.3334
$ws.Range("J122").Value = 1433.3334
$ws.Range("K122").Value = 3799.0002
$ws.Range("L122").Value = 4300.0002
$ws.Range("M122").Value = -1349.0002
$ws.Range("N122").Value = -9200.0002

# Row 132
$ws.Range("H132").Value = 3147.261
$ws.Range("I132").Value = 2779.3333
$ws.Range("J132").Value = 3837.125
$ws.Range("K132").Value = 8337.999899999999
$ws.Range("L132").Value = 11511.375
$ws.Range("M132").Value = -5807.999899999999
$ws.Range("N132").Value = -16571.375

# Row 135
$ws.Range("H135").Value = 40700
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 40700
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 40700
$ws.Range("N135").Value = -50840

# Row 137
$ws.Range("H137").Value = 48745
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 48745
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 48745
$ws.Range("N137").Value = -58945

# Row 139
$ws.Range("H139").Value = 69751.42999999999
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 69751.42999999999
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 69751.42999999999
$ws.Range("N139").Value = -80031.42999999999

# Row 141
$ws.Range("H141").Value = 60315
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 60315
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 60315
$ws.Range("N141").Value = -70675
